$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look numeric,
# so Excel keeps them as literal text instead of coercing to a Double.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.761.70"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.864.59"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("E4").Value = "  -0.80%  "
$ws.Range("D5").Value = "244.52"
$ws.Range("E5").Value = "  -3.81%  "
$ws.Range("D6").Value = "0.676"
$ws.Range("E6").Value = "  -6.92%  "
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("D8").Value = "41.71"
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("E9").Value = "  -4.36%  "
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("D12").Value = "12.85"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "2.138.30"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").Value = "0.710"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "1.865.61"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "34.731.02"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "72.07"
$ws.Range("E18").Value = "  -3.28%  "
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").Value = "242.07"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("E21").Value = "  -4.29%  "
$ws.Range("D22").Value = "4.86"
$ws.Range("E22").Value = "  -4.76%  "
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").Value = "2.46"
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("D25").Value = "2.14"
$ws.Range("E25").Value = "  -14.27%  "
$ws.Range("D26").Value = "162.98"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").Value = "8.31"
$ws.Range("E27").Value = "  -4.01%  "
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  -4.06%  "
$ws.Range("D29").Value = "0.125"
$ws.Range("E29").Value = "  -5.95%  "
$ws.Range("D30").Value = "4.128.45"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "1.69"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").Value = "4.17"
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("D33").Value = "0.0571"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("D35").Value = "4.11"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").Value = "0.824"
$ws.Range("E36").Value = "  -10.63%  "
$ws.Range("E37").Value = "  -20.90%  "
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").Value = "97.25"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").Value = "16.91"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").Value = "0.0210"
$ws.Range("E42").Value = "  -3.71%  "
$ws.Range("E43").Value = "  -5.01%  "
$ws.Range("D44").Value = "1.279.60"
$ws.Range("E44").Value = "  -4.64%  "
$ws.Range("D45").Value = "0.0818"
$ws.Range("E45").Value = "  +10.60%  "
$ws.Range("D46").Value = "2.28"
$ws.Range("E46").Value = "  -6.84%  "
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").Value = "11.69"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("D50").Value = "6.24"
$ws.Range("E50").Value = "  -7.52%  "
$ws.Range("D51").Value = "42.24"
$ws.Range("E51").Value = "  -5.97%  "
